$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.081.10'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.43%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.639.33'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.05%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.01'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.53%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.48'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.505'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.01'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.60%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.250'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.59%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0623'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.71%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0794'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.23%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.760.17'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.27%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.20'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.529'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '62.31'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₃0748'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.00%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.089.60'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.63%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.01'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '190.01'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.57%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.27'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.95%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.69%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '144.27'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.130'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.78%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.59%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.77'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.52%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.74'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.79%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.56%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.53%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0484'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.78%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.93%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.44'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.26%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.50'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.95%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.880'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.43%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.123.14'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.46'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.522'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0155'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.61%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '98.71'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.71%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.786'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.29'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.63%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0₆0114'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.02%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '55.18'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.71%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0521'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.88%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.47'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.415'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.98%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.40%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0927'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.78%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.15'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.33%  '

